$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Return max count for Json entity: update row 5 to the next entry (Sno=2, Name=Banana)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Banana"

# Update the active selection to reflect the cell selected after the edit
$ws.Range("C5").Select()
